$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v3")

# Update the summary formula at the top of the sheet: it used to reference
# the first amortization-schedule date (B10); it now references the second
# one (B11), matching the one-month roll-forward of the payment schedule.
$ws.Range("B1").Formula = "=B11"

# Roll the whole monthly payment-schedule date column forward by one month:
# each row's date becomes the date that used to be in the next row, and a
# new final row is appended one month after the old last row (serial dates,
# precomputed to avoid relying on reading values back through COM).
$ws.Range("B11").Value = 44691
$ws.Range("B12").Value = 44722
$ws.Range("B13").Value = 44752
$ws.Range("B14").Value = 44783
$ws.Range("B15").Value = 44814
$ws.Range("B16").Value = 44844
$ws.Range("B17").Value = 44875
$ws.Range("B18").Value = 44905
$ws.Range("B19").Value = 44936
$ws.Range("B20").Value = 44967
$ws.Range("B21").Value = 44995
$ws.Range("B22").Value = 45026
$ws.Range("B23").Value = 45056
$ws.Range("B24").Value = 45087
$ws.Range("B25").Value = 45117
$ws.Range("B26").Value = 45148
$ws.Range("B27").Value = 45179
$ws.Range("B28").Value = 45209
$ws.Range("B29").Value = 45240
$ws.Range("B30").Value = 45270
$ws.Range("B31").Value = 45301
$ws.Range("B32").Value = 45332
$ws.Range("B33").Value = 45361
$ws.Range("B34").Value = 45392
$ws.Range("B35").Value = 45422
$ws.Range("B36").Value = 45453
$ws.Range("B37").Value = 45483
$ws.Range("B38").Value = 45514
$ws.Range("B39").Value = 45545
$ws.Range("B40").Value = 45575
$ws.Range("B41").Value = 45606
$ws.Range("B42").Value = 45636
$ws.Range("B43").Value = 45667
$ws.Range("B44").Value = 45698
$ws.Range("B45").Value = 45726
$ws.Range("B46").Value = 45757
$ws.Range("B47").Value = 45787
$ws.Range("B48").Value = 45818
$ws.Range("B49").Value = 45848
$ws.Range("B50").Value = 45879
$ws.Range("B51").Value = 45910
$ws.Range("B52").Value = 45940
$ws.Range("B53").Value = 45971
$ws.Range("B54").Value = 46001
$ws.Range("B55").Value = 46032
$ws.Range("B56").Value = 46063
$ws.Range("B57").Value = 46091
$ws.Range("B58").Value = 46122
$ws.Range("B59").Value = 46152
$ws.Range("B60").Value = 46183
$ws.Range("B61").Value = 46213
$ws.Range("B62").Value = 46244
$ws.Range("B63").Value = 46275
$ws.Range("B64").Value = 46305
$ws.Range("B65").Value = 46336
$ws.Range("B66").Value = 46366
$ws.Range("B67").Value = 46397
$ws.Range("B68").Value = 46428
$ws.Range("B69").Value = 46456
$ws.Range("B70").Value = 46487
$ws.Range("B71").Value = 46517
$ws.Range("B72").Value = 46548
$ws.Range("B73").Value = 46578
$ws.Range("B74").Value = 46609
$ws.Range("B75").Value = 46640
$ws.Range("B76").Value = 46670
$ws.Range("B77").Value = 46701
$ws.Range("B78").Value = 46731
$ws.Range("B79").Value = 46762
$ws.Range("B80").Value = 46793
$ws.Range("B81").Value = 46822
$ws.Range("B82").Value = 46853

# Restore the active cell/selection to B1 (it had moved to G12 previously).
$ws.Range("B1").Select()
